$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 227, shifting existing rows 227..302 down to 228..303.
$ws.Rows("227:227").Insert()

# Populate the newly inserted row 227 with the new record.
$ws.Range("A227").Value = 5
$ws.Range("B227").Value = "Macroferia Regional de Talca"
$ws.Range("C227").Value = "Maule"
$ws.Range("D227").Value = 44559
$ws.Range("E227").Value = 7
$ws.Range("F227").Value = 100112043
$ws.Range("G227").Value = "Pepino ensalada"
$ws.Range("H227").Value = "Sin especificar"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 600
$ws.Range("K227").Value = 5000
$ws.Range("L227").Value = 5000
$ws.Range("M227").Value = 5000
$ws.Range("N227").Value = "$/caja 80 unidades"
$ws.Range("O227").Value = "Región del Maule"
$ws.Range("P227").Value = 62
$ws.Range("Q227").Value = 80
$ws.Range("R227").Value = "Hortaliza"
